$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: firstgithubrepo gets a new "User: gauravkhanna3007 / Permission: admin" row
# (previously this row held "secondgithubrepo / Team: secondgithubrepo-teamadmin / ...")
$ws.Range("A4").Value = "github-gk-aks/firstgithubrepo"
$ws.Range("B4").Value = "User: gauravkhanna3007"
$ws.Range("C4").Value = "Permission: admin"

# Row 5: secondgithubrepo / Team: secondgithubrepo-teamadmin (shifted down from old row 4)
$ws.Range("A5").Value = "github-gk-aks/secondgithubrepo"
$ws.Range("B5").Value = "Team: secondgithubrepo-teamadmin"
$ws.Range("C5").Value = "admin: True, maintain: True, push: True, triage: True, pull: True"

# Row 6: secondgithubrepo / Team: secondgithubrepo-teamread (shifted down from old row 5)
$ws.Range("A6").Value = "github-gk-aks/secondgithubrepo"
$ws.Range("B6").Value = "Team: secondgithubrepo-teamread"
$ws.Range("C6").Value = "admin: False, maintain: False, push: False, triage: False, pull: True"

# Row 7: secondgithubrepo gets a new "User: gauravkhanna3007 / Permission: admin" row
$ws.Range("A7").Value = "github-gk-aks/secondgithubrepo"
$ws.Range("B7").Value = "User: gauravkhanna3007"
$ws.Range("C7").Value = "Permission: admin"

# Row 8: thirdgithubrepo / Team: thirdgithubrepo-teamadmin (shifted down from old row 6)
$ws.Range("A8").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B8").Value = "Team: thirdgithubrepo-teamadmin"
$ws.Range("C8").Value = "admin: True, maintain: True, push: True, triage: True, pull: True"

# Row 9: thirdgithubrepo / Team: thirdgithubrepo-teamread (shifted down from old row 7)
$ws.Range("A9").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B9").Value = "Team: thirdgithubrepo-teamread"
$ws.Range("C9").Value = "admin: False, maintain: False, push: False, triage: False, pull: True"

# Row 10: thirdgithubrepo / Team: thirdgithubrepo-teamwrite (shifted down from old row 8)
$ws.Range("A10").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B10").Value = "Team: thirdgithubrepo-teamwrite"
$ws.Range("C10").Value = "admin: False, maintain: False, push: True, triage: True, pull: True"

# Row 11: thirdgithubrepo gets a new "User: arvindsi1973 / Permission: read" row
$ws.Range("A11").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B11").Value = "User: arvindsi1973"
$ws.Range("C11").Value = "Permission: read"

# Row 12: thirdgithubrepo gets a new "User: srinu220kv / Permission: read" row
$ws.Range("A12").Value = "github-gk-aks/thirdgithubrepo"
$ws.Range("B12").Value = "User: srinu220kv"
$ws.Range("C12").Value = "Permission: read"
